# "Add Merch to GP"
#
# The original slide 3 ("IN DEVELOPMENT!") is duplicated. PowerPoint
# inserts the duplicate immediately after the original (position 4,
# fresh SlideID 262) -- we then move that duplicate back to position 3
# so the new merch slide leads and the untouched "IN DEVELOPMENT!"
# slide follows it, matching the new p:sldIdLst order
# (256, 261, 262, 257, 258, 259, 260).

$p = $ppt.ActivePresentation

$original = $p.Slides.Item(3)
$dupRange = $original.Duplicate()
$merch = $dupRange.Item(1)
$merch.MoveTo(3)

# Shape 1: big title ("IN DEVELOPMENT!" -> "WE HAVE MERCH!")
$title = $merch.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "WE HAVE MERCH!"
$title.Left = 410.04930133858267
$title.Top = 30.859616299212597
$title.Width = 404.43638795275587
$title.Height = 65.43283464566929

# Shape 2: body paragraph
$body = $merch.Shapes.Item(2)
$body.TextFrame.TextRange.Text = "We deliver to every country in the world excluding Iran, Cuba, Sudan, North Korea, Syria and Crimea and we have really low prices. This is a great way to support GameProxy. We have stickers, clothing, cases and much more (even a duvet if you really want!)"
$body.Left = 54.4
$body.Top = 96.29244094488189
$body.Width = 1116.7998712598423
$body.Height = 109.05472440944882

# Shape 3: rounded-rectangle call-to-action button
$button = $merch.Shapes.Item(3)
$button.TextFrame.TextRange.Text = "GO TO OUR MERCH STORE"
$button.Left = 447.46306086614175
$button.Top = 302.99528559055113
$button.Width = 329.60820897637797
$button.Height = 61.19763779527559

# Shape 4: caption textbox below the button (text only, no resize/move)
$caption = $merch.Shapes.Item(4)
$caption.TextFrame.TextRange.Text = "It is hosted on Redbubble and we get 15% of the sale."
